$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-style formatting (style index 2, bordered/bold/centered + custom date numfmt) from A269 down through the new date cells A270:A301
$ws.Range("A269").Copy($ws.Range("A270:A301"))

$ws.Range("A270").Value2 = 44344
$ws.Range("B270").Value2 = 1
$ws.Range("C270").Value2 = 1
$ws.Range("D270").Value2 = 16.02307322544464
$ws.Range("A271").Value2 = 44345
$ws.Range("B271").Value2 = 0
$ws.Range("C271").Value2 = 1
$ws.Range("D271").Value2 = 16.02307322544464
$ws.Range("A272").Value2 = 44346
$ws.Range("B272").Value2 = 0
$ws.Range("C272").Value2 = 1
$ws.Range("D272").Value2 = 16.02307322544464
$ws.Range("A273").Value2 = 44347
$ws.Range("B273").Value2 = 0
$ws.Range("C273").Value2 = 1
$ws.Range("D273").Value2 = 16.02307322544464
$ws.Range("A274").Value2 = 44348
$ws.Range("B274").Value2 = 0
$ws.Range("C274").Value2 = 1
$ws.Range("D274").Value2 = 16.02307322544464
$ws.Range("A275").Value2 = 44349
$ws.Range("B275").Value2 = 0
$ws.Range("C275").Value2 = 1
$ws.Range("D275").Value2 = 16.02307322544464
$ws.Range("A276").Value2 = 44350
$ws.Range("B276").Value2 = 1
$ws.Range("C276").Value2 = 2
$ws.Range("D276").Value2 = 32.04614645088928
$ws.Range("A277").Value2 = 44351
$ws.Range("B277").Value2 = 0
$ws.Range("C277").Value2 = 1
$ws.Range("D277").Value2 = 16.02307322544464
$ws.Range("A278").Value2 = 44352
$ws.Range("B278").Value2 = 0
$ws.Range("C278").Value2 = 1
$ws.Range("D278").Value2 = 16.02307322544464
$ws.Range("A279").Value2 = 44353
$ws.Range("B279").Value2 = 0
$ws.Range("C279").Value2 = 1
$ws.Range("D279").Value2 = 16.02307322544464
$ws.Range("A280").Value2 = 44354
$ws.Range("B280").Value2 = 1
$ws.Range("C280").Value2 = 2
$ws.Range("D280").Value2 = 32.04614645088928
$ws.Range("A281").Value2 = 44355
$ws.Range("B281").Value2 = 1
$ws.Range("C281").Value2 = 3
$ws.Range("D281").Value2 = 48.06921967633392
$ws.Range("A282").Value2 = 44356
$ws.Range("B282").Value2 = 1
$ws.Range("C282").Value2 = 4
$ws.Range("D282").Value2 = 64.09229290177856
$ws.Range("A283").Value2 = 44357
$ws.Range("B283").Value2 = 0
$ws.Range("C283").Value2 = 3
$ws.Range("D283").Value2 = 48.06921967633392
$ws.Range("A284").Value2 = 44358
$ws.Range("B284").Value2 = 3
$ws.Range("C284").Value2 = 6
$ws.Range("D284").Value2 = 96.13843935266785
$ws.Range("A285").Value2 = 44359
$ws.Range("B285").Value2 = 0
$ws.Range("C285").Value2 = 6
$ws.Range("D285").Value2 = 96.13843935266785
$ws.Range("A286").Value2 = 44360
$ws.Range("B286").Value2 = 0
$ws.Range("C286").Value2 = 6
$ws.Range("D286").Value2 = 96.13843935266785
$ws.Range("A287").Value2 = 44361
$ws.Range("B287").Value2 = 1
$ws.Range("C287").Value2 = 6
$ws.Range("D287").Value2 = 96.13843935266785
$ws.Range("A288").Value2 = 44362
$ws.Range("B288").Value2 = 0
$ws.Range("C288").Value2 = 5
$ws.Range("D288").Value2 = 80.11536612722321
$ws.Range("A289").Value2 = 44363
$ws.Range("B289").Value2 = 0
$ws.Range("C289").Value2 = 4
$ws.Range("D289").Value2 = 64.09229290177856
$ws.Range("A290").Value2 = 44364
$ws.Range("B290").Value2 = 0
$ws.Range("C290").Value2 = 4
$ws.Range("D290").Value2 = 64.09229290177856
$ws.Range("A291").Value2 = 44365
$ws.Range("B291").Value2 = 0
$ws.Range("C291").Value2 = 1
$ws.Range("D291").Value2 = 16.02307322544464
$ws.Range("A292").Value2 = 44366
$ws.Range("B292").Value2 = 0
$ws.Range("C292").Value2 = 1
$ws.Range("D292").Value2 = 16.02307322544464
$ws.Range("A293").Value2 = 44367
$ws.Range("B293").Value2 = 0
$ws.Range("C293").Value2 = 1
$ws.Range("D293").Value2 = 16.02307322544464
$ws.Range("A294").Value2 = 44368
$ws.Range("B294").Value2 = 0
$ws.Range("C294").Value2 = 0
$ws.Range("D294").Value2 = 0
$ws.Range("A295").Value2 = 44369
$ws.Range("B295").Value2 = 0
$ws.Range("C295").Value2 = 0
$ws.Range("D295").Value2 = 0
$ws.Range("A296").Value2 = 44370
$ws.Range("B296").Value2 = 0
$ws.Range("C296").Value2 = 0
$ws.Range("D296").Value2 = 0
$ws.Range("A297").Value2 = 44371
$ws.Range("B297").Value2 = 0
$ws.Range("C297").Value2 = 0
$ws.Range("D297").Value2 = 0
$ws.Range("A298").Value2 = 44372
$ws.Range("B298").Value2 = 0
$ws.Range("C298").Value2 = 0
$ws.Range("D298").Value2 = 0
$ws.Range("A299").Value2 = 44373
$ws.Range("B299").Value2 = 0
$ws.Range("C299").Value2 = 0
$ws.Range("D299").Value2 = 0
$ws.Range("A300").Value2 = 44374
$ws.Range("B300").Value2 = 0
$ws.Range("C300").Value2 = 0
$ws.Range("D300").Value2 = 0
$ws.Range("A301").Value2 = 44375
$ws.Range("B301").Value2 = 0
$ws.Range("C301").Value2 = 0
$ws.Range("D301").Value2 = 0
